$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ABSM1_RN / M1_RN / CM2_RN / CMN3_RN / CMN4_RN values recomputed after
# removing the < $5 price point (row 12, CMN4_PH = 1.5) from the extrapolation
# calibration input set.

$ws.Range("D3").Value = 119391.2226527308
$ws.Range("E3").Value = -0.000195981089771478
$ws.Range("F3").Value = 0.2005776351199539
$ws.Range("G3").Value = -2.413356533739569
$ws.Range("H3").Value = 34.72329819336328

$ws.Range("D4").Value = 120045.5986622373
$ws.Range("E4").Value = -0.007101267187287857
$ws.Range("F4").Value = 0.2155862234957542
$ws.Range("G4").Value = -1.808078611900467
$ws.Range("H4").Value = 22.88182822685657

$ws.Range("D7").Value = 122018.8542181565
$ws.Range("E7").Value = -0.029298854387986
$ws.Range("F7").Value = 0.2588526672295926
$ws.Range("G7").Value = -1.537140499314168
$ws.Range("H7").Value = 14.98532129978339

$ws.Range("D8").Value = 123618.9583541257
$ws.Range("E8").Value = -0.04432804178162208
$ws.Range("F8").Value = 0.2151841594485501
$ws.Range("G8").Value = -0.8686754940696204
$ws.Range("H8").Value = 7.023883093029358

$ws.Range("D9").Value = 125554.3193767574
$ws.Range("E9").Value = -0.08324576679777931
$ws.Range("F9").Value = 0.3338917264019692
$ws.Range("G9").Value = -1.107842795327797
$ws.Range("H9").Value = 8.13447413961164

$ws.Range("D10").Value = 126840.5479974588
$ws.Range("E10").Value = -0.1255479478825695
$ws.Range("F10").Value = 0.4538497225950301
$ws.Range("G10").Value = -1.71457320976112
$ws.Range("H10").Value = 9.030449326029748

$ws.Range("D11").Value = 129177.22578654
$ws.Range("E11").Value = -0.1365025630979463
$ws.Range("F11").Value = 0.4495110466168772
$ws.Range("G11").Value = -1.450462125962815
$ws.Range("H11").Value = 7.134146175914214

$ws.Range("D14").Value = 118494.7723093694
$ws.Range("E14").Value = 0.03822422267702869
$ws.Range("F14").Value = 0.1286383840365845
$ws.Range("G14").Value = -0.4088091507890652
$ws.Range("H14").Value = 5.505831396975475

$ws.Range("D15").Value = 118552.6968025638
$ws.Range("E15").Value = 0.03664781610329994
$ws.Range("F15").Value = 0.1318328622173739
$ws.Range("G15").Value = -0.01977467477668066
$ws.Range("H15").Value = 5.571520090677694

$ws.Range("D17").Value = 118454.0302185932
$ws.Range("E17").Value = 0.02253468810491821
$ws.Range("F17").Value = 0.130238150300085
$ws.Range("G17").Value = -0.4386338304790469
$ws.Range("H17").Value = 4.694400845141799

$ws.Range("D19").Value = 119418.1356032531
$ws.Range("E19").Value = 0.02264927163941857
$ws.Range("F19").Value = 0.157834041615163
$ws.Range("G19").Value = -0.1618463356375253
$ws.Range("H19").Value = 4.310937060828445
